$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.727457666666666
$ws.Range("H2").Value = 23.182373
$ws.Range("I2").Value = 0.1630271452636819
$ws.Range("J2").Value = 0.1630271452636819
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 36.72315433333333
$ws.Range("N2").Value = 110.169463
$ws.Range("O2").Value = 0.3641786446803374
$ws.Range("P2").Value = 0.3641786446803374
$ws.Range("Q2").Value = 283.7766204972999
$ws.Range("R2").Value = 2553.989584475699
$ws.Range("S2").Value = 0.05937100480823217
$ws.Range("T2").Value = 0.05937100480823217
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.727457666666666
$ws.Range("H3").Value = 23.182373
$ws.Range("I3").Value = 0.1630271452636819
$ws.Range("J3").Value = 0.1630271452636819
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 61.215745
$ws.Range("N3").Value = 183.647235
$ws.Range("O3").Value = 0.6070684137000053
$ws.Range("P3").Value = 0.6070684137000053
$ws.Range("Q3").Value = 473.0420780209616
$ws.Range("R3").Value = 4257.378702188655
$ws.Range("S3").Value = 0.09896863046526372
$ws.Range("T3").Value = 0.0989686304652637
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.727457666666666
$ws.Range("H4").Value = 23.182373
$ws.Range("I4").Value = 0.1630271452636819
$ws.Range("J4").Value = 0.1630271452636819
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.899397666666667
$ws.Range("N4").Value = 8.698193
$ws.Range("O4").Value = 0.02875294161965733
$ws.Range("P4").Value = 0.02875294161965733
$ws.Range("Q4").Value = 22.40497272799878
$ws.Range("R4").Value = 201.644754551989
$ws.Range("S4").Value = 0.004687509990186041
$ws.Range("T4").Value = 0.00468750999018604
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.42779933333333
$ws.Range("H5").Value = 112.283398
$ws.Range("I5").Value = 0.7896189849264272
$ws.Range("J5").Value = 0.7896189849264271
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 36.72315433333333
$ws.Range("N5").Value = 110.169463
$ws.Range("O5").Value = 0.3641786446803374
$ws.Range("P5").Value = 0.3641786446803374
$ws.Range("Q5").Value = 1374.46685127503
$ws.Range("R5").Value = 12370.20166147528
$ws.Range("S5").Value = 0.28756237174437
$ws.Range("T5").Value = 0.28756237174437
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.42779933333333
$ws.Range("H6").Value = 112.283398
$ws.Range("I6").Value = 0.7896189849264272
$ws.Range("J6").Value = 0.7896189849264271
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 61.215745
$ws.Range("N6").Value = 183.647235
$ws.Range("O6").Value = 0.6070684137000053
$ws.Range("P6").Value = 0.6070684137000053
$ws.Range("Q6").Value = 2291.170619900503
$ws.Range("R6").Value = 20620.53557910453
$ws.Range("S6").Value = 0.4793527446066946
$ws.Range("T6").Value = 0.4793527446066945
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.42779933333333
$ws.Range("H7").Value = 112.283398
$ws.Range("I7").Value = 0.7896189849264272
$ws.Range("J7").Value = 0.7896189849264271
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.899397666666667
$ws.Range("N7").Value = 8.698193
$ws.Range("O7").Value = 0.02875294161965733
$ws.Range("P7").Value = 0.02875294161965733
$ws.Range("Q7").Value = 108.5180740555349
$ws.Range("R7").Value = 976.662666499814
$ws.Range("S7").Value = 0.02270386857536264
$ws.Range("T7").Value = 0.02270386857536264
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.244565
$ws.Range("H8").Value = 6.733695
$ws.Range("I8").Value = 0.04735386980989085
$ws.Range("J8").Value = 0.04735386980989083
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 36.72315433333333
$ws.Range("N8").Value = 110.169463
$ws.Range("O8").Value = 0.3641786446803374
$ws.Range("P8").Value = 0.3641786446803374
$ws.Range("Q8").Value = 82.42750690619835
$ws.Range("R8").Value = 741.847562155785
$ws.Range("S8").Value = 0.0172452681277352
$ws.Range("T8").Value = 0.01724526812773519
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.244565
$ws.Range("H9").Value = 6.733695
$ws.Range("I9").Value = 0.04735386980989085
$ws.Range("J9").Value = 0.04735386980989083
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 61.215745
$ws.Range("N9").Value = 183.647235
$ws.Range("O9").Value = 0.6070684137000053
$ws.Range("P9").Value = 0.6070684137000053
$ws.Range("Q9").Value = 137.402718675925
$ws.Range("R9").Value = 1236.624468083325
$ws.Range("S9").Value = 0.02874703862804701
$ws.Range("T9").Value = 0.028747038628047
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.244565
$ws.Range("H10").Value = 6.733695
$ws.Range("I10").Value = 0.04735386980989085
$ws.Range("J10").Value = 0.04735386980989083
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.899397666666667
$ws.Range("N10").Value = 8.698193
$ws.Range("O10").Value = 0.02875294161965733
$ws.Range("P10").Value = 0.02875294161965733
$ws.Range("Q10").Value = 6.507886523681667
$ws.Range("R10").Value = 58.570978713135
$ws.Range("S10").Value = 0.001361563054108645
$ws.Range("T10").Value = 0.001361563054108645
